# typedQNameIRI.xlsx — fix failing atomic tabOTTR tests due to breaking change.
#
# The "data"/"iri" pair in column A (rows 7-9) rotates:
#   A7: "data"  -> 1        (now a numeric literal, was the text "data")
#   A8: 1       -> "iri"    (was the numeric literal, now text "iri")
#   A9: "iri"   -> "data"   (was text "iri", now text "data")
# Existing cell formatting (styles) is left untouched - only the values move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 1
$ws.Range("A8").Value = "iri"
$ws.Range("A9").Value = "data"

# Selection moved from A10 to A8.
$ws.Range("A8").Select()

# Column A got narrower.
$ws.Columns.Item(1).ColumnWidth = 7.6
